# Apply the edits described by the diff:
# - Update error code values in column A for rows 65-70
# - Update the window view (size/position) of the workbook
# - Update the active sheet's scroll position (topLeftCell) and selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update error code values (A65:A70)
$ws.Range("A65").Value = 1000
$ws.Range("A66").Value = 1001
$ws.Range("A67").Value = 1100
$ws.Range("A68").Value = 1101
$ws.Range("A69").Value = 1102
$ws.Range("A70").Value = 1103

# Scroll the window so row 53 is at the top-left, and select A71
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A71").Select()

# Update workbook window position/size
$excel.Left = -37335
$excel.Top = 1433
$excel.Width = 32812
$excel.Height = 11857
